# Auto-generated update of cryptos worksheet values
# Applies the per-cell text changes described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''30.420.75'
$ws.Range("E2").Value = '  -0.93%  '

# Row 3
$ws.Range("D3").Value = '''1.917.24'
$ws.Range("E3").Value = '  +2.05%  '

# Row 4
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
$ws.Range("D5").Value = '''241.69'
$ws.Range("E5").Value = '  +1.74%  '

# Row 6
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  +0.10%  '

# Row 7
$ws.Range("D7").Value = '''0.4702'
$ws.Range("E7").Value = '  -0.90%  '

# Row 8
$ws.Range("D8").Value = '''0.2853'
$ws.Range("E8").Value = '  +0.90%  '

# Row 9
$ws.Range("D9").Value = '''0.06821'
$ws.Range("E9").Value = '  +4.80%  '

# Row 10
$ws.Range("D10").Value = '''107.56'
$ws.Range("E10").Value = '  +12.31%  '

# Row 11
$ws.Range("D11").Value = '''18.27'
$ws.Range("E11").Value = '  -2.04%  '

# Row 12
$ws.Range("D12").Value = '''1.895.49'
$ws.Range("E12").Value = '  +0.93%  '

# Row 13
$ws.Range("D13").Value = '''0.07664'
$ws.Range("E13").Value = '  +1.12%  '

# Row 14
$ws.Range("D14").Value = '''5.205'
$ws.Range("E14").Value = '  +2.43%  '

# Row 15
$ws.Range("D15").Value = '''0.6561'
$ws.Range("E15").Value = '  +1.08%  '

# Row 16
$ws.Range("D16").Value = '''289.76'
$ws.Range("E16").Value = '  -6.10%  '

# Row 17
$ws.Range("D17").Value = '''30.434.39'
$ws.Range("E17").Value = '  -0.90%  '

# Row 18
$ws.Range("D18").Value = '''0.000007630'
$ws.Range("E18").Value = '  +1.29%  '

# Row 19
$ws.Range("D19").Value = '''0.9998'
$ws.Range("E19").Value = '  +0.10%  '

# Row 20
$ws.Range("D20").Value = '''12.94'
$ws.Range("E20").Value = '  -1.06%  '

# Row 21
$ws.Range("D21").Value = '''2.149.96'
$ws.Range("E21").Value = '  +1.27%  '

# Row 22
$ws.Range("D22").Value = '''1.001'
$ws.Range("E22").Value = '  +0.22%  '

# Row 23
$ws.Range("D23").Value = '''5.227'
$ws.Range("E23").Value = '  +1.56%  '

# Row 24
$ws.Range("D24").Value = '''6.198'
$ws.Range("E24").Value = '  +0.63%  '

# Row 25
$ws.Range("D25").Value = '''21.67'
$ws.Range("E25").Value = '  +9.61%  '

# Row 26
$ws.Range("D26").Value = '''168.13'
$ws.Range("E26").Value = '  -0.57%  '

# Row 27
$ws.Range("D27").Value = '''9.294'
$ws.Range("E27").Value = '  +0.45%  '

# Row 28
$ws.Range("D28").Value = '''2.061'
$ws.Range("E28").Value = '  +5.76%  '

# Row 29
$ws.Range("D29").Value = '''0.1074'
$ws.Range("E29").Value = '  +1.45%  '

# Row 30
$ws.Range("D30").Value = '''1.370'
$ws.Range("E30").Value = '  +1.37%  '

# Row 31
$ws.Range("D31").Value = '''4.154'
$ws.Range("E31").Value = '  -0.39%  '

# Row 32
$ws.Range("D32").Value = '''3.966'
$ws.Range("E32").Value = '  +0.50%  '

# Row 33
$ws.Range("D33").Value = '''0.05055'
$ws.Range("E33").Value = '  +0.32%  '

# Row 34
$ws.Range("D34").Value = '''0.7419'
$ws.Range("E34").Value = '  +2.91%  '

# Row 35
$ws.Range("D35").Value = '''1.152'
$ws.Range("E35").Value = '  -1.86%  '

# Row 36
$ws.Range("D36").Value = '''0.02080'
$ws.Range("E36").Value = '  +8.45%  '

# Row 37
$ws.Range("D37").Value = '''2.747'
$ws.Range("E37").Value = '  +1.51%  '

# Row 38
$ws.Range("D38").Value = '''2.692'
$ws.Range("E38").Value = '  -0.26%  '

# Row 39
$ws.Range("D39").Value = '''2.053'
$ws.Range("E39").Value = '  +0.26%  '

# Row 40
$ws.Range("B40").Value = 'Quant'
$ws.Range("C40").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D40").Value = '''108.93'
$ws.Range("E40").Value = '  +1.74%  '

# Row 41
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '''0.8708'
$ws.Range("E41").Value = '  -2.93%  '

# Row 42
$ws.Range("D42").Value = '''5.860'
$ws.Range("E42").Value = '  +4.81%  '

# Row 43
$ws.Range("D43").Value = '''0.9999'
$ws.Range("E43").Value = '  +0.13%  '

# Row 44
$ws.Range("D44").Value = '''0.4236'
$ws.Range("E44").Value = '  +1.00%  '

# Row 45
$ws.Range("D45").Value = '''67.54'
$ws.Range("E45").Value = '  +3.72%  '

# Row 46
$ws.Range("D46").Value = '''50.60'
$ws.Range("E46").Value = '  +18.37%  '

# Row 47
$ws.Range("D47").Value = '''7.172'
$ws.Range("E47").Value = '  -2.13%  '

# Row 48
$ws.Range("D48").Value = '''9.234'
$ws.Range("E48").Value = '  +3.17%  '

# Row 49
$ws.Range("D49").Value = '''0.1212'
$ws.Range("E49").Value = '  -0.29%  '

# Row 50
$ws.Range("D50").Value = '''34.79'
$ws.Range("E50").Value = '  +0.67%  '

# Row 51
$ws.Range("D51").Value = '''0.3903'
$ws.Range("E51").Value = '  +2.57%  '
